# Generate Report for Handback
# Update timestamp cells on the "Overview", "zh-cn" and "de-de" sheets
# to reflect the latest handoff/handback generation times.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-21 09:09:57"

# --- zh-cn sheet: Correspond Handoff/Handback Datetime columns (H, K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-21 09:09:53"
$wsZhCn.Range("K2").Value = "2016-08-21 09:10:17"

# --- de-de sheet: Correspond Handoff/Handback Datetime columns (H, K) ---
# Note: H2 on de-de shares the same original text as Overview!G2
# ("2016-08-21 09:09:15"), so it must be updated too to stay in sync.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-21 09:09:57"
$wsDeDe.Range("K2").Value = "2016-08-21 09:10:24"
